# "combined into one file" -- the framework's compartments/parameters were
# trimmed down (comp_4/comp_5 and par_10..par_14 removed), so the sample
# "Compartments", "Transitions", "Characteristics" and "Parameters" sheets
# that are generated from them shrink accordingly.

$wb = $excel.ActiveWorkbook

# --- Compartments sheet: drop the comp_4 / comp_5 rows (rows 6-7) ---
$wsCompartments = $wb.Worksheets.Item("Compartments")
$wsCompartments.Rows("6:7").Delete()

# --- Transitions sheet: drop the comp_4 / comp_5 column (F:G) and rows (6-7) ---
$wsTransitions = $wb.Worksheets.Item("Transitions")
$wsTransitions.Columns("F:G").Delete()
$wsTransitions.Rows("6:7").Delete()

# --- Characteristics sheet: rows 6 & 7 no longer derive their "components"
# from the now-missing compartments; they fall back to the same chained
# pattern used by the rows beneath them (C<n> = A<n-1>, D<n> empty).
$wsCharacteristics = $wb.Worksheets.Item("Characteristics")
$wsCharacteristics.Range("C6").Formula = "=A5"
$wsCharacteristics.Range("D6").ClearContents()
$wsCharacteristics.Range("C7").Formula = "=A6"
$wsCharacteristics.Range("D7").ClearContents()

# --- Parameters sheet: drop the par_10 .. par_14 rows (12-16) ---
$wsParameters = $wb.Worksheets.Item("Parameters")
$wsParameters.Rows("12:16").Delete()
